$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74

$ws.Cells.Item($row, 1).Value = 43578
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 2).Value = 3
$ws.Cells.Item($row, 3).Value = 63
$ws.Cells.Item($row, 4).Value = 77
$ws.Cells.Item($row, 5).Value = 129
$ws.Cells.Item($row, 6).Value = 68
$ws.Cells.Item($row, 7).Value = 17
$ws.Cells.Item($row, 8).Value = 6
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 75
$ws.Cells.Item($row, 11).Value = 37
$ws.Cells.Item($row, 12).Value = 38
$ws.Cells.Item($row, 13).Value = 26
$ws.Cells.Item($row, 14).Value = 12
$ws.Cells.Item($row, 15).Value = 4

$ws.Range("H76").Select()
